$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '96.608.81'
$ws.Range("E2").Value = '  -0.71%  '

$ws.Range("D3").Value = '3.658.52'
$ws.Range("E3").Value = '  -1.39%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '2.64'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +37.97%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.999'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.11%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '227.19'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.63%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '643.34'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.57%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.425'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.24%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.14'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.70%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.999'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.01%  '

$ws.Range("D11").Value = '3.657.10'
$ws.Range("E11").Value = '  -1.37%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '48.76'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +9.65%  '

$ws.Range("E13").Value = '  +0.97%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000294'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -8.12%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.65'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.30%  '

$ws.Range("D16").Value = '4.326.15'
$ws.Range("E16").Value = '  -1.72%  '

$ws.Range("D17").Value = '95.904.18'
$ws.Range("E17").Value = '  -1.21%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '21.54'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +14.90%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.89'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.34%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.25'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +8.39%  '

$ws.Range("D21").Value = '3.640.94'
$ws.Range("E21").Value = '  -1.79%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.538'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.43%  '

$ws.Range("E23").Value = '  +34.66%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '518.78'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.60%  '

$ws.Range("E25").Value = '  -4.51%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '122.48'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +19.41%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000203'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -7.06%  '

$ws.Range("E28").Value = '  -1.02%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '12.96'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.04%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '13.31'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.45%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.01'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.18%  '

$ws.Range("E32").Value = '  +0.07%  '

$ws.Range("E33").Value = '  -3.77%  '

$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '33.26'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.98%  '

$ws.Range("B35").Value = 'PolygonEcosystemToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.625'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.99%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.10%  '

$ws.Range("E37").Value = '  -4.31%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '607.05'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.54%  '

$ws.Range("E39").Value = '  -0.02%  '

$ws.Range("E40").Value = '  -4.20%  '

$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '43.36'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.54%  '

$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '7.15'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.02%  '

$ws.Range("E43").Value = '  +0.30%  '

$ws.Range("B44").Value = 'Kaspa'
$ws.Range("C44").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.161'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.20%  '

$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0500'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.39%  '

$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.963'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.06%  '

$ws.Range("B47").Value = 'ImmutableX'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.97'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.30%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.30'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.36%  '

$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '227.85'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +10.23%  '

$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.83'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.49%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.56'
$ws.Range("D51").Style = "Normal"
